# Update the multiplication problems in the two-digit-mul worksheet.
# Each old expression is unique within the document, so a simple
# Find/Replace (wrap = none needed, whole document scope) is safe.

$d = $word.ActiveDocument

$replacements = @(
    @("13×98=", "90×74="),
    @("47×24=", "72×79="),
    @("71×60=", "17×15="),
    @("30×72=", "51×34="),
    @("15×31=", "75×27="),
    @("20×24=", "76×49="),
    @("65×68=", "52×59="),
    @("56×48=", "88×96="),
    @("82×79=", "90×87="),
    @("35×92=", "70×68="),
    @("95×80=", "74×76="),
    @("81×34=", "18×93="),
    @("40×79=", "35×27="),
    @("49×31=", "37×88="),
    @("76×14=", "53×94="),
    @("12×14=", "11×37="),
    @("14×27=", "47×93="),
    @("74×38=", "16×32="),
    @("37×63=", "21×84="),
    @("71×52=", "44×36="),
    @("91×40=", "96×14="),
    @("27×97=", "82×61="),
    @("95×82=", "68×86="),
    @("74×26=", "17×77="),
    @("70×34=", "78×36=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
